# Legacy GSC export data update:
# The rolling daily coverage window has advanced by one day, so the oldest
# date row (2025-11-02) is dropped from the "Chart" sheet. Deleting the
# entire row shifts every subsequent row up by one, which matches the
# target data (and naturally causes Excel to renumber/clean up the shared
# string table and shrink the sheet's used range to A1:D87 on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date (2025-11-02); remove it and shift everything up.
$ws.Rows.Item(2).Delete()
